# Applies a cyclic rotation of the "event" data in rows 18-21 of the
# Artfynd sheet: the content that used to live in row 19 moves up to
# row 18, row 20's content moves to row 19, row 21's content moves to
# row 20, and row 18's original content wraps around into row 21.
#
# Only the cells that actually change value (per the target diff) are
# touched; cells identical across all four rows (D, I, P, S, T, U, V,
# W, Y, AA, AD, AE, AG, AT, AW, AX, AY) are left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 18 (becomes old row 19's data) ----
$ws.Range("A18").Value = 130979946
$ws.Range("B18").Value = 57884
$ws.Range("E18").Value = 100109
$ws.Range("F18").Value = "Tretåig hackspett"
$ws.Range("G18").Value = "Picoides tridactylus"
$ws.Range("H18").Value = "(Linnaeus, 1758)"
$ws.Range("K18").Value = ""
$ws.Range("L18").Value = ""
$ws.Range("M18").Value = "färska spår"
$ws.Range("N18").Value = ""
$ws.Range("Q18").Value = 590605
$ws.Range("R18").Value = 6963364
$ws.Range("Z18").Value = "09:47"
$ws.Range("AB18").Value = "09:47"
$ws.Range("AC18").Value = "färska ringhack på gran"

# ---- Row 19 (becomes old row 20's data) ----
$ws.Range("A19").Value = 130979899
$ws.Range("Q19").Value = 590850
$ws.Range("R19").Value = 6963133
$ws.Range("Z19").Value = "13:16"
$ws.Range("AB19").Value = "13:16"

# ---- Row 20 (becomes old row 21's data) ----
$ws.Range("A20").Value = 130979897
$ws.Range("B20").Value = 80348
$ws.Range("E20").Value = 6458
$ws.Range("F20").Value = "Lunglav"
$ws.Range("G20").Value = "Lobaria pulmonaria"
$ws.Range("H20").Value = "(L.) Hoffm."
$ws.Range("K20").ClearContents()
$ws.Range("L20").ClearContents()
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("Q20").Value = 590726
$ws.Range("R20").Value = 6963153
$ws.Range("Z20").Value = "13:24"
$ws.Range("AB20").Value = "13:24"
$ws.Range("AC20").ClearContents()

# ---- Row 21 (becomes old row 18's data) ----
$ws.Range("A21").Value = 130979947
$ws.Range("B21").Value = 91808
$ws.Range("E21").Value = 1202
$ws.Range("F21").Value = "Ullticka"
$ws.Range("G21").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H21").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q21").Value = 590591
$ws.Range("R21").Value = 6963354
$ws.Range("Z21").Value = "09:45"
$ws.Range("AB21").Value = "09:45"
